# Automatic update of files.
# Appends a new species observation record as row 35 on the "Artfynd"
# sheet, matching the structure of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

$ws.Cells.Item($row, 1).Value = 111702865
$ws.Cells.Item($row, 2).Value = 56543
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
$ws.Cells.Item($row, 4).Value = "NT"
$ws.Cells.Item($row, 5).Value = 103021
$ws.Cells.Item($row, 6).Value = "Talltita"
$ws.Cells.Item($row, 7).Value = "Poecile montanus"
$ws.Cells.Item($row, 8).Value = "(Conrad von Baldenstein, 1827)"
# Stored as text in the source data ("2"), not a number - force text entry
# with a leading apostrophe so it is not auto-converted to a numeric value.
$ws.Cells.Item($row, 9).Value = "'2"

# Age/stage and sex are present but blank for this record (still present
# as empty text cells in the source data).
$ws.Cells.Item($row, 11).Value = "'"
$ws.Cells.Item($row, 12).Value = "'"
$ws.Cells.Item($row, 13).Value = "födosökande"
$ws.Cells.Item($row, 14).Value = "'"

$ws.Cells.Item($row, 16).Value = "Börtingtjärnen, Nabbnäs, Ly lm"
$ws.Cells.Item($row, 17).Value = 557511
$ws.Cells.Item($row, 18).Value = 7314093
$ws.Cells.Item($row, 19).Value = 127
$ws.Cells.Item($row, 20).Value = "Västerbotten"
$ws.Cells.Item($row, 21).Value = "Sorsele"
$ws.Cells.Item($row, 22).Value = "Lycksele lappmark"
$ws.Cells.Item($row, 23).Value = "Sorsele"

# Dates/times are stored as plain text in this sheet, not real date
# values - force text entry for the date-looking strings so Excel does
# not reinterpret them as date serial numbers.
$ws.Cells.Item($row, 25).Value = "'2023-08-26"
$ws.Cells.Item($row, 26).Value = "14:05"
$ws.Cells.Item($row, 27).Value = "'2023-08-26"
$ws.Cells.Item($row, 28).Value = "14:05"

$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false

# Determination year is blank for this record (empty text cell).
$ws.Cells.Item($row, 46).Value = "'"

$ws.Cells.Item($row, 49).Value = "Erik Owusu-Ansah"
$ws.Cells.Item($row, 50).Value = "Erik Owusu-Ansah"

# Project name is blank for this record (empty text cell).
$ws.Cells.Item($row, 51).Value = "'"

Write-Output "Row 35 written"
